# "add a new holy stone"
#
# The item table on Sheet1 (backed by Excel Table "表2", range A3:AC134)
# gets one new data row inserted right before the current row 92
# (i.e. directly after the existing "符文-艾尔" / 22302018 row), for a
# new "holy stone" item:
#   Id=22302019, ~Name=符文-艾德 (via the table's lookup formula),
#   CdGroup=4, CdTime=15, HolyWord="holyman"
# Every following row shifts down by one (old row 92 -> 93, ..., old
# row 134 -> 135), and the table/used-range grow by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a fresh blank row at sheet row 92; this pushes the old row 92
# ("木质修理锤", 22302030) and everything below it down by one row, and
# extends the sheet's dimension/used range automatically.
$ws.Rows.Item(92).Insert()

# The new row should look like its neighbours (same cell style as the
# rest of the data rows): clone formatting from the row right below it
# (which is the old row 92, now shifted to row 93).
$ws.Range("A93:AC93").Copy()
$ws.Range("A92:AC92").PasteSpecial(-4122)  # xlPasteFormats

# Column B ("~Name") in this table uses a slightly different look for
# freshly-added rows elsewhere in the sheet (e.g. row 116); copy that
# cell's format onto the new B92 so it matches.
$ws.Range("B116").Copy()
$ws.Range("B92").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Fill in the new row's data.
$ws.Range("A92").Value = 22302019
$ws.Range("B92").Formula = '=LOOKUP(表2[[#This Row],[Id]],[1]其他!$A:$A,[1]其他!$B:$B)'
$ws.Range("C92").Value = 4
$ws.Range("D92").Value = 15
$ws.Range("Z92").Value = "holyman"

# Grow the Excel Table (ListObject) so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:AC135"))
